$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: merge the four runs "Lesson 1" / "3" / " - " / "Recursion"
#    (all sharing identical rPr) into a single run "Lesson 13 - Recursion".
#    Find & Replace naturally collapses the matched text into one run while
#    carrying over the original (shared) run formatting.
# ---------------------------------------------------------------------------
$oldTitle = "Lesson 1" + [char]0x0033 + [char]0x0020 + [char]0x2013 + [char]0x0020 + "Recursion"
$newTitle = "Lesson 13" + [char]0x0020 + [char]0x2013 + [char]0x0020 + "Recursion"
$d.Content.Find.Execute($oldTitle, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newTitle, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "May result in a StackOverFlow or Heap Storage Exhaustion error" needs to
#    become three separate runs (no formatting change, just a run split) so
#    that "Stack Overflow" is isolated as its own run, also fixing the
#    casing/spacing of "StackOverFlow" -> "Stack Overflow".
#    A plain Find/Replace (or Range.Text=) would just leave it as a single
#    merged run, so InsertXML (which inserts literal runs verbatim, without
#    adding any rPr) is used on the exact sub-range instead.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*StackOverFlow*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $oldText = "May result in a StackOverFlow or Heap Storage Exhaustion error"
    $start = $target.Range.Start
    $idx = $target.Range.Text.IndexOf($oldText)
    $rangeToReplace = $d.Range($start + $idx, $start + $idx + $oldText.Length)

    $openXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">May result in a </w:t></w:r><w:r><w:t>Stack Overflow</w:t></w:r><w:r><w:t xml:space="preserve"> or Heap Storage Exhaustion error</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $rangeToReplace.InsertXML($openXml)
}
